$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.920.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.887.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4103"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.53"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07981"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9914"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.844.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.911"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.064"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.963.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.379"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.214"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.102.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.119"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.411"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9782"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09357"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.414"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.602"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06055"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02232"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.262"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.176"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5775"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1819"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.262"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.282"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5480"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.909"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07006"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.12%  "
